# Maestro.xlsx update:
#  - Insert a new article row ("Budín marmolado Buon Natale") at row 7 of the
#    "Artículos" sheet, shifting the existing rows (previously 7..62) down to 8..63.
#  - Fix a wrong image filename that was already present for the "Café torrado
#    molido" article (now row 39): it pointed to 7790150006153.png, should be
#    7790150006351.png (matching the article's own barcode).
#  - A couple of boolean flags (Pesable / TieneVencimiento columns) get flipped
#    from FALSE to TRUE for the existing rows 5 and 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artículos")

# --- Insert the new row and shift everything below it down ---
$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value  = 7798094222820
$ws.Cells.Item(7, 2).Value  = "Budín"
$ws.Cells.Item(7, 3).Value  = "marmolado"
$ws.Cells.Item(7, 4).Value  = "clásico"
$ws.Cells.Item(7, 5).Value  = "Buon Natale"
$ws.Cells.Item(7, 6).Value  = 170
$ws.Cells.Item(7, 7).Value  = "gr."
$ws.Cells.Item(7, 8).Value  = "bolsa"
$ws.Cells.Item(7, 9).Value  = "Budines"
$ws.Cells.Item(7, 10).Value = "Argentina"
$ws.Cells.Item(7, 11).Value = 12
$ws.Cells.Item(7, 12).Value = $false
$ws.Cells.Item(7, 13).Value = $true
$ws.Cells.Item(7, 14).Value = "C:\VentaSoft\Imágenes de artículos\7798094222820.png"
$ws.Cells.Item(7, 15).Value = $true
$ws.Cells.Item(7, 16).Value = $true

# --- Fix the mismatched image path for the "Café" article (now row 39) ---
$ws.Cells.Item(39, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790150006351.png"

# --- Flip a couple of boolean flags on the pre-existing rows 5 and 6 ---
$ws.Cells.Item(5, 16).Value = $true
$ws.Cells.Item(6, 13).Value = $true
$ws.Cells.Item(6, 15).Value = $true

# --- Column E ("Marca") widened slightly to fit the new brand name ("Buon Natale") ---
$ws.Columns.Item(5).ColumnWidth = 11
